# Update the header date
$d = $word.ActiveDocument
[void]$d.Content.Find.Execute("2025-08-12 Tuesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-08-13 Wednesday", 2)

# Update each answer cell in the table directly by (row, column) address so that
# duplicate values elsewhere in the table are not accidentally disturbed.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "91÷6=15, 1"
$t.Cell(1,2).Range.Text = "36÷3=12, 0"
$t.Cell(1,3).Range.Text = "37÷8=4, 5"
$t.Cell(1,4).Range.Text = "87÷5=17, 2"
$t.Cell(1,5).Range.Text = "90÷5=18, 0"

$t.Cell(5,1).Range.Text = "27÷6=4, 3"
$t.Cell(5,2).Range.Text = "11÷3=3, 2"
$t.Cell(5,3).Range.Text = "65÷4=16, 1"
$t.Cell(5,4).Range.Text = "27÷6=4, 3"
$t.Cell(5,5).Range.Text = "47÷9=5, 2"

$t.Cell(9,1).Range.Text = "80÷6=13, 2"
$t.Cell(9,2).Range.Text = "21÷6=3, 3"
$t.Cell(9,3).Range.Text = "70÷5=14, 0"
$t.Cell(9,4).Range.Text = "47÷4=11, 3"
$t.Cell(9,5).Range.Text = "17÷4=4, 1"

$t.Cell(13,1).Range.Text = "22÷6=3, 4"
$t.Cell(13,2).Range.Text = "11÷8=1, 3"
$t.Cell(13,3).Range.Text = "98÷7=14, 0"
$t.Cell(13,4).Range.Text = "46÷6=7, 4"
$t.Cell(13,5).Range.Text = "67÷7=9, 4"

$t.Cell(17,1).Range.Text = "43÷9=4, 7"
$t.Cell(17,2).Range.Text = "65÷5=13, 0"
$t.Cell(17,3).Range.Text = "36÷4=9, 0"
$t.Cell(17,4).Range.Text = "68÷8=8, 4"
$t.Cell(17,5).Range.Text = "87÷5=17, 2"

Write-Output "ok"
